# Refresh the cryptocurrency price ("D") and 1h volume change ("E")
# columns with the latest scraped figures.
#
# The sheet stores these as plain text, preserving formatting such as
# trailing zeros, "thousand.thousand.unit"-style prices, and the padded
# "  +x.xx%  " change strings. Excel's Range.Value setter auto-detects
# plain decimal-looking strings (e.g. "246.36") and silently converts
# them to numbers, which would corrupt the intended text (dropping
# trailing zeros, switching to scientific notation, etc.). To avoid
# that, cells whose new value looks like a plain number are first
# switched to Text number format ("@") so the literal text is kept.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '30.569.55' }
    @{ Cell = 'E2'; Value = '  -0.60%  ' }
    @{ Cell = 'D3'; Value = '1.885.22' }
    @{ Cell = 'E3'; Value = '  -0.27%  ' }
    @{ Cell = 'E4'; Value = '  -0.08%  ' }
    @{ Cell = 'D5'; Value = '246.36' }
    @{ Cell = 'E5'; Value = '  -0.78%  ' }
    @{ Cell = 'E6'; Value = '  -0.02%  ' }
    @{ Cell = 'D7'; Value = '0.4742' }
    @{ Cell = 'E7'; Value = '  +0.06%  ' }
    @{ Cell = 'D8'; Value = '0.2895' }
    @{ Cell = 'E8'; Value = '  -1.10%  ' }
    @{ Cell = 'D9'; Value = '0.06544' }
    @{ Cell = 'E9'; Value = '  +0.11%  ' }
    @{ Cell = 'D10'; Value = '22.38' }
    @{ Cell = 'E10'; Value = '  +1.11%  ' }
    @{ Cell = 'D11'; Value = '0.7749' }
    @{ Cell = 'E11'; Value = '  +4.84%  ' }
    @{ Cell = 'D12'; Value = '100.98' }
    @{ Cell = 'E12'; Value = '  +4.17%  ' }
    @{ Cell = 'D14'; Value = '1.884.78' }
    @{ Cell = 'E14'; Value = '  -0.26%  ' }
    @{ Cell = 'D15'; Value = '5.264' }
    @{ Cell = 'E15'; Value = '  +0.30%  ' }
    @{ Cell = 'D16'; Value = '285.23' }
    @{ Cell = 'E16'; Value = '  -0.04%  ' }
    @{ Cell = 'D17'; Value = '30.554.95' }
    @{ Cell = 'E18'; Value = '  -0.34%  ' }
    @{ Cell = 'D19'; Value = '0.000007540' }
    @{ Cell = 'E19'; Value = '  -0.13%  ' }
    @{ Cell = 'E20'; Value = '  -0.03%  ' }
    @{ Cell = 'D21'; Value = '2.131.08' }
    @{ Cell = 'E21'; Value = '  -0.21%  ' }
    @{ Cell = 'D22'; Value = '5.365' }
    @{ Cell = 'E22'; Value = '  +0.70%  ' }
    @{ Cell = 'E23'; Value = '  +0.15%  ' }
    @{ Cell = 'D24'; Value = '6.462' }
    @{ Cell = 'E24'; Value = '  +3.21%  ' }
    @{ Cell = 'D25'; Value = '9.178' }
    @{ Cell = 'E25'; Value = '  -0.56%  ' }
    @{ Cell = 'D26'; Value = '163.33' }
    @{ Cell = 'E26'; Value = '  -1.01%  ' }
    @{ Cell = 'D27'; Value = '19.15' }
    @{ Cell = 'E27'; Value = '  +0.69%  ' }
    @{ Cell = 'D28'; Value = '1.917' }
    @{ Cell = 'E28'; Value = '  -0.08%  ' }
    @{ Cell = 'D29'; Value = '1.336' }
    @{ Cell = 'E29'; Value = '  -0.34%  ' }
    @{ Cell = 'D30'; Value = '0.09710' }
    @{ Cell = 'E30'; Value = '  -0.44%  ' }
    @{ Cell = 'D31'; Value = '1.502' }
    @{ Cell = 'E31'; Value = '  +0.38%  ' }
    @{ Cell = 'D32'; Value = '4.259' }
    @{ Cell = 'E32'; Value = '  -0.92%  ' }
    @{ Cell = 'D33'; Value = '4.194' }
    @{ Cell = 'E33'; Value = '  +0.22%  ' }
    @{ Cell = 'D34'; Value = '0.04850' }
    @{ Cell = 'E34'; Value = '  -0.31%  ' }
    @{ Cell = 'D35'; Value = '1.132' }
    @{ Cell = 'E35'; Value = '  +0.54%  ' }
    @{ Cell = 'D36'; Value = '0.6993' }
    @{ Cell = 'E36'; Value = '  +0.28%  ' }
    @{ Cell = 'D37'; Value = '2.765' }
    @{ Cell = 'E37'; Value = '  +1.53%  ' }
    @{ Cell = 'D38'; Value = '0.01919' }
    @{ Cell = 'E38'; Value = '  +1.38%  ' }
    @{ Cell = 'E39'; Value = '  +3.29%  ' }
    @{ Cell = 'D40'; Value = '76.04' }
    @{ Cell = 'E40'; Value = '  -0.18%  ' }
    @{ Cell = 'D41'; Value = '6.301' }
    @{ Cell = 'E41'; Value = '  -0.45%  ' }
    @{ Cell = 'E42'; Value = '  -0.12%  ' }
    @{ Cell = 'D43'; Value = '0.4260' }
    @{ Cell = 'E43'; Value = '  -0.49%  ' }
    @{ Cell = 'E44'; Value = '  -0.08%  ' }
    @{ Cell = 'D45'; Value = '0.8335' }
    @{ Cell = 'E45'; Value = '  -0.19%  ' }
    @{ Cell = 'D46'; Value = '101.65' }
    @{ Cell = 'D47'; Value = '9.910' }
    @{ Cell = 'E47'; Value = '  +3.91%  ' }
    @{ Cell = 'D48'; Value = '7.031' }
    @{ Cell = 'E48'; Value = '  -0.31%  ' }
    @{ Cell = 'D49'; Value = '35.24' }
    @{ Cell = 'E49'; Value = '  -1.08%  ' }
    @{ Cell = 'D50'; Value = '895.00' }
    @{ Cell = 'E50'; Value = '  -2.18%  ' }
    @{ Cell = 'D51'; Value = '0.05773' }
    @{ Cell = 'E51'; Value = '  +0.23%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.Value.Trim() -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a plain number -> force Text format so Excel keeps
        # the exact original string instead of re-parsing it as a number.
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
